$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.476.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.510.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.37%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.510.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.90%  "

$ws.Range("E11").Value = "  -0.79%  "

$ws.Range("E12").Value = "  -3.78%  "

$ws.Range("E13").Value = "  -1.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.969.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.413.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.66%  "

$ws.Range("E16").Value = "  -6.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.502.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.31%  "

$ws.Range("E19").Value = "  -6.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "347.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.36%  "

$ws.Range("E23").Value = "  -4.64%  "

$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.60%  "

$ws.Range("E27").Value = "  -7.35%  "

$ws.Range("E28").Value = "  -4.94%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("E30").Value = "  -6.02%  "

$ws.Range("E31").Value = "  -1.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "471.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.16%  "

$ws.Range("E33").Value = "  -0.78%  "

$ws.Range("E34").Value = "  -3.36%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("E36").Value = "  +2.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.31%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("E41").Value = "  -3.31%  "

$ws.Range("E42").Value = "  -3.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.07%  "

$ws.Range("E44").Value = "  -13.65%  "

$ws.Range("E45").Value = "  -9.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.528"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.89%  "

$ws.Range("E50").Value = "  -5.37%  "

$ws.Range("E51").Value = "  -2.21%  "
